$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "I0" (I1) and "IF" (J1) with the same formatting as
# --- the existing header cells (bold font, thin border, center/top align).
# Copy the format of the neighboring header cell (H1) onto I1 and J1 first,
# then set their text so the shared style index is reused rather than a new
# one being created.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-34: new I (I0) / J (IF) columns.
$data = @(
    @(2, 1, 6),
    @(3, 8, 9),
    @(4, 6, 8),
    @(5, 4, 6),
    @(6, 1, 4),
    @(7, 2, 6),
    @(8, 1, 5),
    @(9, 8, 8),
    @(10, 1, 4),
    @(11, 1, 5),
    @(12, 1, 5),
    @(13, 1, 6),
    @(14, 1, 6),
    @(15, 1, 7),
    @(16, 1, 5),
    @(17, 1, 5),
    @(18, 1, 5),
    @(19, 1, 5),
    @(20, 1, 5),
    @(21, 1, 6),
    @(22, 1, 3),
    @(23, 1, 5),
    @(24, 1, 6),
    @(25, 1, 5),
    @(26, 6, 6),
    @(27, 1, 5),
    @(28, 1, 2),
    @(29, 1, 4),
    @(30, 1, 5),
    @(31, 1, 4),
    @(32, 1, 4),
    @(33, 1, 3),
    @(34, 1, 2)
)

foreach ($row in $data) {
    $r = $row[0]
    $i = $row[1]
    $j = $row[2]
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
}
